$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$measureName = "add_ideal_loads"

# Columns: B=Argument name, C=Display name, D=Type, E=Unit, F=Description, G=Required, H=Default
$rows = @(
    @("heat_recovery_method", "Heat recovery method", "String", "", 'The method by which heat recovery is done. One of: "None", "Sensible" or "Enthalpy".', "True", "none"),
    @("latent_efficiency", "Latent efficiency", "Float", "", "The efficiency of heat recovery for the latent method.", "True", "0.65"),
    @("sensible_efficiency", "Sensible efficiency", "Float", "", "The efficiency of heat recovery for the sensible method.", "True", "0.7"),
    @("ach_per_hour", "Air changes per hour", "Float", "h^-1", "Number of air changes per hour through mechanical ventilation.", "True", "1.0"),
    @("nfa_gfa_ratio", "Ratio of NFA over GFA", "Float", "", "Ratio of NFA over GFA.", "True", "1.0"),
    @("floor_height_ratio", "Ratio of conditioned floor height over total floor height", "Float", "", "Ratio of conditioned floor height over total floor height.", "True", "1.0"),
    @("hvac_schedule", "HVAC schedule name", "String", "", "(Export only) Name of the HVAC schedule", "False", ""),
    @("is_custom_hvac", "Is custom HVAC schedule", "Bool", "", "(Export only) Flag whether the HVAC schedule is custom or from the standard library.", "False", ""),
    @("hvac_sched_weekday", "HVAC schedule for weekdays", "String", "", "HVAC schedule definition for weekdays.", "True", ""),
    @("hvac_sched_saturday", "HVAC schedule for saturday", "String", "", "HVAC schedule definition for saturday.", "True", ""),
    @("hvac_sched_sunday", "HVAC schedule for sunday", "String", "", "HVAC schedule definition for sunday.", "True", ""),
    @("hvac_sched_holiday", "HVAC schedule for holiday", "String", "", "HVAC schedule definition for holidays.", "False", ""),
    @("holidays", "Holiday definition", "String", "", "Definition of holidays", "False", "")
)

# Find first empty row right after the current used range (row 114 -> 115)
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $measureName
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    if ($row[3] -ne "") {
        $ws.Cells.Item($r, 5).Value = $row[3]
    }
    $ws.Cells.Item($r, 6).Value = $row[4]

    # Columns G (Required) and H (Default) hold text that looks like a
    # boolean/number ("True"/"False"/"0.65"/...). A direct .Value assignment
    # gets auto-coerced to a Boolean/Number cell by the engine's input
    # parser, but the source workbook stores these as plain shared strings.
    # Writing them as a formula that yields the literal text, then
    # paste-special-values over itself, keeps the cell text-typed.
    $ws.Cells.Item($r, 7).Formula = '="' + $row[5] + '"'
    if ($row[6] -ne "") {
        $ws.Cells.Item($r, 8).Formula = '="' + $row[6] + '"'
    }

    $r = $r + 1
}

$lastRow = $r - 1
$ghRange = $ws.Range($ws.Cells.Item($startRow, 7), $ws.Cells.Item($lastRow, 8))
[void]$ghRange.Copy()
[void]$ghRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Mirror the new "one past the last used row" active cell / selection the
# author's Excel session ended up with after appending the rows.
[void]$ws.Cells.Item($lastRow + 1, 1).Select()
